$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells keep their original plain-text representation
# (several values look like numbers, e.g. "21.90" or "0.990", and would
# otherwise be auto-converted to numeric values losing the trailing zero).
$priceCells = @("D2","D3","D5","D8","D9","D12","D13","D15","D16","D17","D18","D19","D22","D23","D25","D26","D27","D32","D33","D37","D39","D42","D43","D46","D47","D48","D49")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.750.53'
$ws.Range("E2").Value = '  -2.60%  '
$ws.Range("D3").Value = '1.565.53'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '206.32'
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("E6").Value = '  -2.36%  '
$ws.Range("D8").Value = '21.90'
$ws.Range("E8").Value = '  -0.85%  '
$ws.Range("D9").Value = '0.247'
$ws.Range("E9").Value = '  -0.82%  '
$ws.Range("E10").Value = '  -1.31%  '
$ws.Range("E11").Value = '  -0.32%  '
$ws.Range("D12").Value = '1.788.28'
$ws.Range("E12").Value = '  -0.05%  '
$ws.Range("D13").Value = '1.570.71'
$ws.Range("E13").Value = '  +0.35%  '
$ws.Range("E14").Value = '  -2.58%  '
$ws.Range("D15").Value = '0.514'
$ws.Range("E15").Value = '  -0.84%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '26.794.27'
$ws.Range("E16").Value = '  -2.41%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = '61.45'
$ws.Range("E17").Value = '  -3.43%  '
$ws.Range("D18").Value = '214.86'
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("D19").Value = '7.37'
$ws.Range("E19").Value = '  +1.70%  '
$ws.Range("E20").Value = '  -1.75%  '
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("D22").Value = '4.09'
$ws.Range("E22").Value = '  -0.31%  '
$ws.Range("D23").Value = '9.30'
$ws.Range("E23").Value = '  -2.74%  '
$ws.Range("E24").Value = '  -1.36%  '
$ws.Range("D25").Value = '152.57'
$ws.Range("E25").Value = '  -0.53%  '
$ws.Range("D26").Value = '6.73'
$ws.Range("E26").Value = '  +0.85%  '
$ws.Range("D27").Value = '14.92'
$ws.Range("E27").Value = '  -0.36%  '
$ws.Range("E28").Value = '  +0.14%  '
$ws.Range("E29").Value = '  -1.38%  '
$ws.Range("E31").Value = '  -1.20%  '
$ws.Range("D32").Value = '3.15'
$ws.Range("E32").Value = '  -1.54%  '
$ws.Range("D33").Value = '1.389.46'
$ws.Range("E33").Value = '  +1.03%  '
$ws.Range("E34").Value = '  -1.34%  '
$ws.Range("E35").Value = '  +0.33%  '
$ws.Range("E36").Value = '  -0.84%  '
$ws.Range("D37").Value = '0.931'
$ws.Range("E37").Value = '  -2.72%  '
$ws.Range("E38").Value = '  -2.52%  '
$ws.Range("D39").Value = '0.527'
$ws.Range("E39").Value = '  -0.63%  '
$ws.Range("E40").Value = '  -0.60%  '
$ws.Range("E41").Value = '  +0.15%  '
$ws.Range("D42").Value = '0.990'
$ws.Range("E42").Value = '  +1.23%  '
$ws.Range("D43").Value = '1.79'
$ws.Range("E43").Value = '  -0.52%  '
$ws.Range("E44").Value = '  +1.74%  '
$ws.Range("E45").Value = '  +0.78%  '
$ws.Range("D46").Value = '63.29'
$ws.Range("E46").Value = '  -1.29%  '
$ws.Range("D47").Value = '1.701.26'
$ws.Range("E47").Value = '  +0.13%  '
$ws.Range("D48").Value = '85.65'
$ws.Range("E48").Value = '  +0.19%  '
$ws.Range("D49").Value = '0.0₇0991'
$ws.Range("E49").Value = '  -1.18%  '
$ws.Range("E50").Value = '  -0.87%  '
